$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Move Robot2 to location (11, 8) and remove the toolkit."
$ws.Range("E2").Value = $true

$ws.Range("A3").Value = "Move Robot26 to location (4, 4) and remove the liquid spill."

$ws.Range("A4").Value = "Move Robot42 to location (9, 1) and remove the large debris."

$ws.Range("A5").Value = "Move Robot50 to location (7, 11) and remove the dust."
$ws.Range("B5").Value = $true
$ws.Range("C5").Value = $true

$ws.Range("A6").Value = "Move Robot41 to location (6, 12) and remove the grass."

$ws.Range("A7").Value = "Move Robot50 to location (3, 1) and remove the small debris."
$ws.Range("B7").Value = $true
$ws.Range("C7").Value = $true

$ws.Range("A8").Value = "Move Robot13 to location (1, 4) and remove the vehicle."

$ws.Range("A9").Value = "Move Robot13 to location (11, 1) and remove the construction materials."

$ws.Range("A10").Value = "Move Robot14 to location (2, 10) and remove the tree branches."

$ws.Range("A11").Value = "Move Robot15 to location (8, 6) and remove the screws."
